# Weekly update: insert a new price record as the first data row (row 175)
# for "Hortaliza, Terminal La Palmera de La Serena - Cebollín".
#
# The new observation is inserted above the existing history, pushing all
# subsequent rows (old 175..235) down by one (new 176..236), and the sheet
# dimension grows from A1:R235 to A1:R236.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 175; this shifts rows 175:235 down to 176:236
# and extends the used range/dimension automatically.
$ws.Rows("175:175").Insert()

# Populate the new row 175 with this week's data.
$ws.Range("A175").Value = 8
$ws.Range("B175").Value = "Terminal La Palmera de La Serena"
$ws.Range("C175").Value = "Coquimbo"
$ws.Range("D175").Value = 44809
$ws.Range("E175").Value = 4
$ws.Range("F175").Value = 100112037
$ws.Range("G175").Value = "Cebollín"
$ws.Range("H175").Value = "Sin especificar"
$ws.Range("I175").Value = "Primera"
$ws.Range("J175").Value = 1400
$ws.Range("K175").Value = 1400
$ws.Range("L175").Value = 1600
$ws.Range("M175").Value = 1500
$ws.Range("N175").Value = "$/paquete 6 unidades"
$ws.Range("O175").Value = "Provincia del Elquí"
$ws.Range("P175").Value = 250
$ws.Range("Q175").Value = 6
$ws.Range("R175").Value = "Hortaliza"
